$wb = $excel.ActiveWorkbook

# 1. Add the new "SPRINT 7" sheet, placed after the last existing sheet ("SPRINT 6")
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "SPRINT 7"

# Reference sheet to copy date-cell formatting from (keeps same style index as other sprint sheets)
$refSheet = $wb.Worksheets.Item("SPRINT 6")

# 2. Header / total formulas
$ws.Range("A1").Formula = "=SUM(D3:D14)"

$ws.Range("B2").Value = "Date"
$ws.Range("C2").Value = "Quoi"
$ws.Range("D2").Value = "Temps (h)"

# 3. Data rows
$ws.Range("C3").Value = "Sprint review + comments"
$ws.Range("D3").Value = 1

$ws.Range("C4").Value = "Planification"
$ws.Range("D4").Value = 1

$ws.Range("C5").Value = "Créer rapport avec API"
$ws.Range("D5").Value = 6

$ws.Range("C6").Value = "Créer rapport avec API"
$ws.Range("D6").Value = 6

$ws.Range("C7").Value = "Planification, git"
$ws.Range("D7").Value = 0.5

# 4. Date column (B3:B16), formatted like the other sprint sheets (copy format then set values)
$refSheet.Range("B3:B16").Copy()
$ws.Range("B3:B16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B3").Value = 44015
$ws.Range("B4").Value = 44015
$ws.Range("B5").Value = 44021
$ws.Range("B6").Value = 44025
$ws.Range("B7").Value = 44025

# 5. Footer total
$ws.Range("D17").Formula = "=SUM(D3:D16)"

$ws.Range("I26").Select()
